$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H32").Value = 1740.25
$ws.Range("I32").Value = 850.5
$ws.Range("J32").Value = 2036.8334
$ws.Range("K32").Value = 850.5
$ws.Range("L32").Value = 2036.8334
$ws.Range("M32").Value = -524.5
$ws.Range("N32").Value = -2688.8334
$ws.Range("H116").Value = 4466183.5
$ws.Range("I116").Value = 5767970.5
$ws.Range("J116").Value = 2914.5715
$ws.Range("K116").Value = 5767970.5
$ws.Range("L116").Value = 2914.5715
$ws.Range("M116").Value = -5764528.5
$ws.Range("N116").Value = -9798.5715
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("M121").ClearContents()
$ws.Range("H132").Value = 197401.14
$ws.Range("I132").Value = 213816.98
$ws.Range("J132").Value = 80438.25
$ws.Range("K132").Value = 641450.9400000001
$ws.Range("L132").Value = 241314.75
$ws.Range("M132").Value = -638920.9400000001
$ws.Range("N132").Value = -246374.75
$ws.Range("H133").Value = 19671.818
$ws.Range("J133").Value = 19671.818
$ws.Range("L133").Value = 19671.818
$ws.Range("N133").Value = -29791.818
$ws.Range("H135").Value = 1163.6459
$ws.Range("I135").Value = 1056.2559
$ws.Range("J135").Value = 2087.2
$ws.Range("K135").Value = 9506.303100000001
$ws.Range("L135").Value = 18784.8
$ws.Range("M135").Value = -6971.303100000001
$ws.Range("N135").Value = -23854.8
$ws.Range("H137").Value = 43479410
$ws.Range("I137").Value = 45455624
$ws.Range("K137").Value = 136366872
$ws.Range("M137").Value = -136364322
$ws.Range("H138").Value = 4123514
$ws.Range("I138").Value = 1070356.9
$ws.Range("J138").Value = 6413382
$ws.Range("K138").Value = 3211070.7
$ws.Range("L138").Value = 19240146
$ws.Range("M138").Value = -3205930.7
$ws.Range("N138").Value = -19250426

$ws = $wb.Worksheets.Item(2)
$ws.Range("H45").Value = 775.7692
$ws.Range("I45").Value = 698.63635
$ws.Range("J45").Value = 1200
$ws.Range("K45").Value = 698.63635
$ws.Range("L45").Value = 1200
$ws.Range("M45").Value = -321.63635
$ws.Range("N45").Value = -1954
$ws.Range("H74").Value = 4768.816
$ws.Range("I74").Value = 1428.4642
$ws.Range("J74").Value = 14121.8
$ws.Range("K74").Value = 1428.4642
$ws.Range("L74").Value = 14121.8
$ws.Range("M74").Value = -554.4641999999999
$ws.Range("N74").Value = -15869.8
$ws.Range("H77").Value = 4768.816
$ws.Range("I77").Value = 1428.4642
$ws.Range("J77").Value = 14121.8
$ws.Range("K77").Value = 7142.321
$ws.Range("L77").Value = 70609
$ws.Range("M77").Value = -2774.321
$ws.Range("N77").Value = -79345
$ws.Range("H122").Value = 2466.6191
$ws.Range("I122").Value = 2163.4443
$ws.Range("J122").Value = 2694
$ws.Range("K122").Value = 6490.3329
$ws.Range("L122").Value = 8082
$ws.Range("M122").Value = -4040.3329
$ws.Range("N122").Value = -12982
$ws.Range("H132").Value = 2070.6667
$ws.Range("I132").Value = 1573.2122
$ws.Range("J132").Value = 4806.6665
$ws.Range("K132").Value = 4719.6366
$ws.Range("L132").Value = 14419.9995
$ws.Range("M132").Value = -2189.6366
$ws.Range("N132").Value = -19479.9995
$ws.Range("H133").Value = 49800
$ws.Range("J133").Value = 49800
$ws.Range("L133").Value = 49800
$ws.Range("N133").Value = -54860
$ws.Range("H139").Value = 45119.168
$ws.Range("J139").Value = 45119.168
$ws.Range("L139").Value = 45119.168
$ws.Range("N139").Value = -55399.168

$ws = $wb.Worksheets.Item(3)
$ws.Range("H107").Value = 667.5714
$ws.Range("I107").Value = 693.3333
$ws.Range("K107").Value = 693.3333
$ws.Range("M107").Value = 1226.6667
$ws.Range("H133").Value = 78393.336
$ws.Range("J133").Value = 78393.336
$ws.Range("L133").Value = 78393.336
$ws.Range("N133").Value = -88513.336
$ws.Range("H134").Value = 15627133
$ws.Range("I134").Value = 22728652
$ws.Range("K134").Value = 68185956
$ws.Range("M134").Value = -68183421

$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 2394.4707
$ws.Range("I31").Value = 1245.7778
$ws.Range("K31").Value = 1245.7778
$ws.Range("M31").Value = -950.7778000000001
$ws.Range("H34").Value = 2394.4707
$ws.Range("I34").Value = 1245.7778
$ws.Range("K34").Value = 1245.7778
$ws.Range("M34").Value = -1043.7778
$ws.Range("H122").Value = 2127.3125
$ws.Range("I122").Value = 1169.1111
$ws.Range("J122").Value = 3359.2856
$ws.Range("K122").Value = 3507.3333
$ws.Range("L122").Value = 10077.8568
$ws.Range("M122").Value = -1057.3333
$ws.Range("N122").Value = -14977.8568
$ws.Range("H134").Value = 1799.3971
$ws.Range("I134").Value = 1125.5103
$ws.Range("J134").Value = 3537.3157
$ws.Range("K134").Value = 3376.5309
$ws.Range("L134").Value = 10611.9471
$ws.Range("M134").Value = -841.5308999999997
$ws.Range("N134").Value = -15681.9471

$ws = $wb.Worksheets.Item(5)
$ws.Range("H32").Value = 4333.3335
$ws.Range("J32").Value = 4333.3335
$ws.Range("L32").Value = 13000.0005
$ws.Range("N32").Value = -13566.0005
$ws.Range("H92").Value = 800
$ws.Range("I92").Value = 701
$ws.Range("J92").Value = 866
$ws.Range("K92").Value = 2103
$ws.Range("L92").Value = 2598
$ws.Range("M92").Value = -855
$ws.Range("N92").Value = -5094
$ws.Range("H113").Value = 14706796
$ws.Range("I113").Value = 631.5333000000001
$ws.Range("J113").Value = 26316926
$ws.Range("K113").Value = 1894.5999
$ws.Range("L113").Value = 78950778
$ws.Range("M113").Value = 275.4000999999998
$ws.Range("N113").Value = -78955118
$ws.Range("H122").Value = 638.25
$ws.Range("I122").Value = 275.72726
$ws.Range("J122").Value = 1081.3334
$ws.Range("K122").Value = 2481.54534
$ws.Range("L122").Value = 9732.000599999999
$ws.Range("M122").Value = -31.54534000000012
$ws.Range("N122").Value = -14632.0006

$ws = $wb.Worksheets.Item(6)
$ws.Range("H122").Value = 1112262.4
$ws.Range("I122").Value = 1588317.8
$ws.Range("J122").Value = 1466.6666
$ws.Range("K122").Value = 4764953.4
$ws.Range("L122").Value = 4399.9998
$ws.Range("M122").Value = -4762503.4
$ws.Range("N122").Value = -9299.9998
$ws.Range("H137").Value = 54850
$ws.Range("J137").Value = 54850
$ws.Range("L137").Value = 54850
$ws.Range("N137").Value = -65050
$ws.Range("H138").Value = 64133.332
$ws.Range("J138").Value = 64133.332
$ws.Range("L138").Value = 64133.332
$ws.Range("N138").Value = -74413.33199999999
$ws.Range("H139").Value = 36163
$ws.Range("J139").Value = 36163
$ws.Range("L139").Value = 36163
$ws.Range("N139").Value = -46443

$ws = $wb.Worksheets.Item(7)
$ws.Range("H40").Value = 2280.139
$ws.Range("I40").Value = 1379
$ws.Range("J40").Value = 3406.5625
$ws.Range("K40").Value = 1379
$ws.Range("L40").Value = 3406.5625
$ws.Range("M40").Value = -1243
$ws.Range("N40").Value = -3678.5625
$ws.Range("H61").Value = 7174.6787
$ws.Range("I61").Value = 7095.091
$ws.Range("J61").Value = 7466.5
$ws.Range("K61").Value = 7095.091
$ws.Range("L61").Value = 7466.5
$ws.Range("M61").Value = -6893.091
$ws.Range("N61").Value = -7870.5
$ws.Range("H82").Value = 1105
$ws.Range("I82").Value = 900
$ws.Range("J82").Value = 1222.1428
$ws.Range("K82").Value = 900
$ws.Range("L82").Value = 1222.1428
$ws.Range("M82").Value = -539
$ws.Range("N82").Value = -1944.1428
$ws.Range("H85").Value = 1105
$ws.Range("I85").Value = 900
$ws.Range("J85").Value = 1222.1428
$ws.Range("K85").Value = 900
$ws.Range("L85").Value = 1222.1428
$ws.Range("M85").Value = 348
$ws.Range("N85").Value = -3718.1428
$ws.Range("H113").Value = 7174.6787
$ws.Range("I113").Value = 7095.091
$ws.Range("J113").Value = 7466.5
$ws.Range("K113").Value = 7095.091
$ws.Range("L113").Value = 7466.5
$ws.Range("M113").Value = -4925.091
$ws.Range("N113").Value = -11806.5
$ws.Range("H136").Value = 4240.268
$ws.Range("I136").Value = 2570.303
$ws.Range("K136").Value = 7710.909
$ws.Range("M136").Value = -5160.909

$ws = $wb.Worksheets.Item(8)
$ws.Range("H62").Value = 16686217
$ws.Range("I62").Value = 25026950
$ws.Range("J62").Value = 4750
$ws.Range("K62").Value = 25026950
$ws.Range("L62").Value = 4750
$ws.Range("M62").Value = -25026326
$ws.Range("N62").Value = -5998
$ws.Range("H65").Value = 16686217
$ws.Range("I65").Value = 25026950
$ws.Range("J65").Value = 4750
$ws.Range("K65").Value = 125134750
$ws.Range("L65").Value = 23750
$ws.Range("M65").Value = -125131630
$ws.Range("N65").Value = -29990
$ws.Range("H128").Value = 42245.832
$ws.Range("J128").Value = 42245.832
$ws.Range("L128").Value = 42245.832
$ws.Range("N128").Value = -52205.832
$ws.Range("H136").Value = 9553802
$ws.Range("I136").Value = 11145374
$ws.Range("J136").Value = 4368.4
$ws.Range("K136").Value = 33436122
$ws.Range("L136").Value = 13105.2
$ws.Range("M136").Value = -33433572
$ws.Range("N136").Value = -18205.2
